# Update ticker lists on Sheet1 and drop the now-unused trailing rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B (Buying Opportunity) and C (support Zone) for rows 2-15.
$colB = @(
  "NSE:DCAL",
  "NSE:GENUSPOWER",
  "NSE:GREENPLY",
  "NSE:HINDWAREAP",
  "NSE:HIRECT",
  "NSE:ICIL",
  "NSE:LATENTVIEW",
  "NSE:MATRIMONY",
  "NSE:RITCO",
  "NSE:ROSSARI",
  "NSE:RVNL",
  "",
  "",
  ""
)

$colC = @(
  "NSE:ANMOL",
  "NSE:BAJFINANCE",
  "NSE:DBOL",
  "NSE:EMAMIREAL",
  "NSE:HDFCLOWVOL",
  "NSE:HINDZINC",
  "NSE:HITECHCORP",
  "NSE:HNGSNGBEES",
  "NSE:LICI",
  "NSE:MANKIND",
  "NSE:MONTECARLO",
  "NSE:NEWGEN",
  "NSE:PNBHOUSING",
  "NSE:RKEC"
)

# New values for column E (Short buildup) for rows 2-8; rows 9-15 stay blank.
$colE = @(
  "NSE:ADANIENT",
  "NSE:DALBHARAT",
  "NSE:GMRINFRA",
  "NSE:INDHOTEL",
  "NSE:JUBLFOOD",
  "NSE:MARICO",
  "NSE:MUTHOOTFIN"
)

for ($i = 0; $i -lt $colB.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 2).Value = $colB[$i]
  $ws.Cells.Item($row, 3).Value = $colC[$i]
}

for ($i = 0; $i -lt $colE.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 5).Value = $colE[$i]
}

# Column F (FII ENTERING) is cleared for every data row - the content moved
# to column E and rows 16-19 (which only had column C populated) are removed.
$ws.Range("F2:F15").ClearContents()

# Rows 16-19 no longer exist in the updated sheet.
$ws.Range("A16:A19").EntireRow.Delete()
